# Append the next day's GSC export row ("2026-01-04") to the bottom of
# the "Chart" sheet, right after the existing last row (90, "2026-01-03").
# The "Table" sheet (Issue/Validation/Pages header) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Chart" sheet

$row = 91

# Column A holds dates stored as plain text (e.g. "2026-01-03"), not real
# date serials. Assigning a date-shaped string straight to .Value would be
# auto-parsed into a date number, so force the cell to Text first, write
# the value, then drop the temporary Text format again so the cell ends
# up styled exactly like its neighbours (default/general style).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2026-01-04"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = 0.0
$ws.Cells.Item($row, 3).Value = 27.0
